$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.469.97'
$ws.Range("E2").Value = '  +0.41%  '
$ws.Range("D3").Value = '1.808.77'
$ws.Range("E3").Value = '  +0.26%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '225.41'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.589'
$ws.Range("E6").Value = '  +2.57%  '
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '38.47'
$ws.Range("E8").Value = '  +6.68%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.288'
$ws.Range("E9").Value = '  -4.34%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0673'
$ws.Range("E10").Value = '  -2.84%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0974'
$ws.Range("E11").Value = '  +0.81%  '
$ws.Range("D12").Value = '2.071.52'
$ws.Range("E12").Value = '  +0.31%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '11.14'
$ws.Range("E13").Value = '  -4.86%  '
$ws.Range("D14").Value = '1.808.09'
$ws.Range("E14").Value = '  +0.21%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.630'
$ws.Range("E15").Value = '  -2.28%  '
$ws.Range("D16").Value = '34.463.14'
$ws.Range("E16").Value = '  +0.44%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.42'
$ws.Range("E17").Value = '  -1.55%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '68.16'
$ws.Range("E18").Value = '  -1.27%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '242.95'
$ws.Range("E19").Value = '  -1.00%  '
$ws.Range("D20").Value = '0.0₃0772'
$ws.Range("E20").Value = '  -2.77%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.19'
$ws.Range("E21").Value = '  -3.04%  '
$ws.Range("E22").Value = '  -0.10%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.11'
$ws.Range("E23").Value = '  -1.80%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.20'
$ws.Range("E24").Value = '  +3.33%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '170.26'
$ws.Range("E25").Value = '  -1.14%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.73'
$ws.Range("E26").Value = '  -2.73%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.66'
$ws.Range("E27").Value = '  +4.68%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.120'
$ws.Range("E28").Value = '  +1.57%  '
$ws.Range("E29").Value = '  -0.04%  '
$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.23'
$ws.Range("E30").Value = '  -1.52%  '
$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.78'
$ws.Range("E31").Value = '  -1.81%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0516'
$ws.Range("E32").Value = '  -2.96%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.85'
$ws.Range("E33").Value = '  -4.36%  '
$ws.Range("E34").Value = '  -0.38%  '
$ws.Range("D35").Value = '1.354.38'
$ws.Range("E35").Value = '  -2.74%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.640'
$ws.Range("E36").Value = '  -4.90%  '
$ws.Range("E37").Value = '  -0.59%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0187'
$ws.Range("E38").Value = '  -1.64%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.33'
$ws.Range("E39").Value = '  -5.08%  '
$ws.Range("E40").Value = '  +1.43%  '
$ws.Range("B41").Value = 'WEMIXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.21'
$ws.Range("E41").Value = '  -1.61%  '
$ws.Range("B42").Value = 'ARBITRUM'
$ws.Range("C42").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.951'
$ws.Range("E42").Value = '  -1.29%  '
$ws.Range("B43").Value = 'Aave'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '81.75'
$ws.Range("E43").Value = '  -0.05%  '
$ws.Range("E44").Value = '  -0.95%  '
$ws.Range("E45").Value = '  +1.02%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0511'
$ws.Range("E46").Value = '  +1.70%  '
$ws.Range("D47").Value = '1.972.78'
$ws.Range("E47").Value = '  +0.35%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.77'
$ws.Range("E48").Value = '  -4.31%  '
$ws.Range("E49").Value = '  -0.10%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '102.31'
$ws.Range("E50").Value = '  -2.32%  '
$ws.Range("E51").Value = '  -4.89%  '
